$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 46
$ws.Range("I6").Value = 46
$ws.Range("K6").Value = 138
$ws.Range("M6").Value = -26

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2465
$ws.Range("I32").Value = 2475
$ws.Range("J32").Value = 2445
$ws.Range("K32").Value = 2475
$ws.Range("L32").Value = 2445
$ws.Range("M32").Value = -2149
$ws.Range("N32").Value = -3097

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 30033.334
$ws.Range("I82").Value = 22500
$ws.Range("J82").Value = 45100
$ws.Range("K82").Value = 67500
$ws.Range("L82").Value = 135300
$ws.Range("M82").Value = -67094
$ws.Range("N82").Value = -136112

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 30033.334
$ws.Range("I85").Value = 22500
$ws.Range("J85").Value = 45100
$ws.Range("K85").Value = 67500
$ws.Range("L85").Value = 135300
$ws.Range("M85").Value = -66096
$ws.Range("N85").Value = -138108

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1375
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2785.7144
$ws.Range("I138").Value = 1500
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 4500
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = 640
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4351.154
$ws.Range("I32").Value = 4187.727
$ws.Range("K32").Value = 4187.727
$ws.Range("M32").Value = -3900.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2962.25
$ws.Range("I45").Value = 850
$ws.Range("K45").Value = 850
$ws.Range("M45").Value = -473

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 624.6667
$ws.Range("I74").Value = 624.6667
$ws.Range("K74").Value = 624.6667
$ws.Range("M74").Value = 249.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 624.6667
$ws.Range("I77").Value = 624.6667
$ws.Range("K77").Value = 3123.3335
$ws.Range("M77").Value = 1244.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 50110
$ws.Range("J80").Value = 50110
$ws.Range("L80").Value = 50110
$ws.Range("N80").Value = -52106

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 50110
$ws.Range("J83").Value = 50110
$ws.Range("L83").Value = 150330
$ws.Range("N83").Value = -160314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 20000000
$ws.Range("J112").Value = 20000000
$ws.Range("L112").Value = 20000000
$ws.Range("N112").Value = -20002954

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 52332.25
$ws.Range("J124").Value = 52332.25
$ws.Range("L124").Value = 52332.25
$ws.Range("N124").Value = -62152.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 3000
$ws.Range("J21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("N21").Value = -3470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3896.7144
$ws.Range("I31").Value = 3455.8
$ws.Range("J31").Value = 4999
$ws.Range("K31").Value = 3455.8
$ws.Range("L31").Value = 4999
$ws.Range("M31").Value = -3160.8
$ws.Range("N31").Value = -5589

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3896.7144
$ws.Range("I34").Value = 3455.8
$ws.Range("J34").Value = 4999
$ws.Range("K34").Value = 3455.8
$ws.Range("L34").Value = 4999
$ws.Range("M34").Value = -3253.8
$ws.Range("N34").Value = -5403

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1349.5
$ws.Range("I99").Value = 1199
$ws.Range("K99").Value = 1199
$ws.Range("M99").Value = 299

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1349.5
$ws.Range("I126").Value = 1199
$ws.Range("K126").Value = 3597
$ws.Range("M126").Value = -1127

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 793
$ws.Range("I5").Value = 680
$ws.Range("J5").Value = 811.8333
$ws.Range("K5").Value = 2040
$ws.Range("L5").Value = 2435.4999
$ws.Range("M5").Value = -1928
$ws.Range("N5").Value = -2659.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 81.666664
$ws.Range("I33").Value = 83.333336
$ws.Range("J33").Value = 80
$ws.Range("K33").Value = 500.000016
$ws.Range("L33").Value = 480
$ws.Range("M33").Value = -217.000016
$ws.Range("N33").Value = -1046

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 342.5
$ws.Range("I44").Value = 294.57144
$ws.Range("K44").Value = 883.71432
$ws.Range("M44").Value = -485.71432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 220
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4935
$ws.Range("I129").Value = 2497.5
$ws.Range("J129").Value = 5747.5
$ws.Range("K129").Value = 7492.5
$ws.Range("L129").Value = 17242.5
$ws.Range("M129").Value = -2492.5
$ws.Range("N129").Value = -27242.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1629.091
$ws.Range("J131").Value = 1886.25
$ws.Range("L131").Value = 5658.75
$ws.Range("N131").Value = -15738.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 793
$ws.Range("I135").Value = 680
$ws.Range("J135").Value = 811.8333
$ws.Range("K135").Value = 6120
$ws.Range("L135").Value = 7306.4997
$ws.Range("M135").Value = -3585
$ws.Range("N135").Value = -12376.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3299.2856
$ws.Range("I43").Value = 4000
$ws.Range("K43").Value = 4000
$ws.Range("M43").Value = -3849

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1297
$ws.Range("I122").Value = 1382.5714
$ws.Range("J122").Value = 997.5
$ws.Range("K122").Value = 4147.7142
$ws.Range("L122").Value = 2992.5
$ws.Range("M122").Value = -1697.7142
$ws.Range("N122").Value = -7892.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 13334
$ws.Range("I2").Value = 20000
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = -19888
$ws.Range("N2").Value = -226

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3406.5
$ws.Range("I40").Value = 3406.5
$ws.Range("K40").Value = 3406.5
$ws.Range("M40").Value = -3270.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2002
$ws.Range("J46").Value = 2002
$ws.Range("L46").Value = 2002
$ws.Range("N46").Value = -2378

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4331.6665
$ws.Range("I93").Value = 4331.6665
$ws.Range("K93").Value = 4331.6665
$ws.Range("M93").Value = -3083.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2677.6667
$ws.Range("I122").Value = 1516.5
$ws.Range("K122").Value = 4549.5
$ws.Range("M122").Value = -2099.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11985.4
$ws.Range("I136").Value = 7731.125
$ws.Range("K136").Value = 23193.375
$ws.Range("M136").Value = -20643.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1818649.5
$ws.Range("I14").Value = 2352956.5
$ws.Range("J14").Value = 2006
$ws.Range("K14").Value = 2352956.5
$ws.Range("L14").Value = 2006
$ws.Range("M14").Value = -2352788.5
$ws.Range("N14").Value = -2342

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4153.875
$ws.Range("I126").Value = 2621.8333
$ws.Range("K126").Value = 7865.499899999999
$ws.Range("M126").Value = -5395.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2584.5715
$ws.Range("I132").Value = 2584.5715
$ws.Range("K132").Value = 7753.7145
$ws.Range("M132").Value = -5223.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2445.3333
$ws.Range("I136").Value = 2188.5
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 6565.5
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -4015.5
$ws.Range("N136").Value = -18600
